$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 68
$ws1.Range("F8").Value = 105
$ws1.Range("F9").Value = 8683
$ws1.Range("F11").Value = 329
$ws1.Range("F13").Value = 972
$ws1.Range("F14").Value = 108
$ws1.Range("F17").Value = 233
$ws1.Range("F18").Value = 249
$ws1.Range("F21").Value = 1013

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 68
$ws4.Range("F10").Value = 105
$ws4.Range("F11").Value = 8683
$ws4.Range("F13").Value = 329
$ws4.Range("F15").Value = 972
$ws4.Range("F16").Value = 108
$ws4.Range("F19").Value = 233
$ws4.Range("F20").Value = 249
$ws4.Range("F23").Value = 1013
